# MainIcon.xlsx edit:
#   "add some new equips. remove some useless module"
#
# The only real content change in the sheet is the removal of the
# "MainIcon6" / 奇术 ("Spellcraft") entry, which lived in row 10
# (Id=9). Removing the worksheet row:
#   - shifts every following row up by one (Id values / data untouched),
#   - shrinks the used range from A1:L23 to A1:L22,
#   - shrinks table "表1" and its AutoFilter from A1:L23 to A1:L22,
#   - drops the now-unreferenced shared strings
#     ("奇术", "查看我的奇术", "MainIcon6") from sharedStrings.xml,
#   - and leaves the selection on the row that used to be row 11
#     (now row 10), matching Excel's default "select the row below
#     the deleted one" behaviour after a row delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete worksheet row 10 (Id=9, Name="奇术"/MainIcon6) and shift
# everything below it up by one row.
$ws.Rows.Item(10).Delete() | Out-Null

# Match the post-edit selection recorded in the workbook.
$ws.Range("A10:XFD10").Select() | Out-Null
